# A new weekly price-report row was inserted at row 54 ("Fruta / hortaliza,
# semanal"), pushing the previously existing rows 54-101 down to 55-102.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 54, shifting rows 54:101 -> 55:102.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
$ws.Cells.Item(54, 1).Value = 3
$ws.Cells.Item(54, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value = 45040
$ws.Cells.Item(54, 5).Value = 5
$ws.Cells.Item(54, 6).Value = 100112035
$ws.Cells.Item(54, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 40
$ws.Cells.Item(54, 11).Value = 12000
$ws.Cells.Item(54, 12).Value = 12000
$ws.Cells.Item(54, 13).Value = 12000
$ws.Cells.Item(54, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(54, 16).Value = 1200
$ws.Cells.Item(54, 17).Value = 10
$ws.Cells.Item(54, 18).Value = "Hortaliza"
